$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new text value, and whether it must be forced to
# plain text (the "Price" column holds values such as "308.01" or "1.000" that
# Excel would otherwise auto-convert to numbers, losing the trailing zeros / dots).
$updates = @(
    @{ Cell = "D2"; Value = "23.574.37"; ForceText = $true }
    @{ Cell = "D3"; Value = "1.641.40"; ForceText = $true }
    @{ Cell = "E3"; Value = "  +2.42%  "; ForceText = $false }
    @{ Cell = "D4"; Value = "0.9996"; ForceText = $true }
    @{ Cell = "E4"; Value = "  +0.10%  "; ForceText = $false }
    @{ Cell = "D5"; Value = "308.01"; ForceText = $true }
    @{ Cell = "E5"; Value = "  +1.55%  "; ForceText = $false }
    @{ Cell = "D6"; Value = "1.000"; ForceText = $true }
    @{ Cell = "E6"; Value = "  +0.12%  "; ForceText = $false }
    @{ Cell = "D7"; Value = "0.3761"; ForceText = $true }
    @{ Cell = "E7"; Value = "  -0.50%  "; ForceText = $false }
    @{ Cell = "E8"; Value = "  +1.18%  "; ForceText = $false }
    @{ Cell = "D9"; Value = "0.3667"; ForceText = $true }
    @{ Cell = "E9"; Value = "  +1.52%  "; ForceText = $false }
    @{ Cell = "E10"; Value = "  +1.05%  "; ForceText = $false }
    @{ Cell = "D11"; Value = "0.08203"; ForceText = $true }
    @{ Cell = "E11"; Value = "  +0.98%  "; ForceText = $false }
    @{ Cell = "D12"; Value = "0.9995"; ForceText = $true }
    @{ Cell = "E12"; Value = "  +0.09%  "; ForceText = $false }
    @{ Cell = "D13"; Value = "23.05"; ForceText = $true }
    @{ Cell = "E13"; Value = "  +1.65%  "; ForceText = $false }
    @{ Cell = "D14"; Value = "6.672"; ForceText = $true }
    @{ Cell = "E14"; Value = "  +0.99%  "; ForceText = $false }
    @{ Cell = "D15"; Value = "0.00001287"; ForceText = $true }
    @{ Cell = "E15"; Value = "  +2.73%  "; ForceText = $false }
    @{ Cell = "D16"; Value = "7.435"; ForceText = $true }
    @{ Cell = "E16"; Value = "  +0.43%  "; ForceText = $false }
    @{ Cell = "D17"; Value = "1.646.08"; ForceText = $true }
    @{ Cell = "E17"; Value = "  +2.58%  "; ForceText = $false }
    @{ Cell = "D18"; Value = "94.97"; ForceText = $true }
    @{ Cell = "E18"; Value = "  +1.59%  "; ForceText = $false }
    @{ Cell = "D19"; Value = "0.06923"; ForceText = $true }
    @{ Cell = "E19"; Value = "  +0.88%  "; ForceText = $false }
    @{ Cell = "D20"; Value = "18.30"; ForceText = $true }
    @{ Cell = "E20"; Value = "  +1.32%  "; ForceText = $false }
    @{ Cell = "D21"; Value = "6.576"; ForceText = $true }
    @{ Cell = "E21"; Value = "  +0.68%  "; ForceText = $false }
    @{ Cell = "D22"; Value = "0.9988"; ForceText = $true }
    @{ Cell = "D23"; Value = "23.574.23"; ForceText = $true }
    @{ Cell = "E23"; Value = "  +1.53%  "; ForceText = $false }
    @{ Cell = "E24"; Value = "  -0.74%  "; ForceText = $false }
    @{ Cell = "D25"; Value = "3.079"; ForceText = $true }
    @{ Cell = "E25"; Value = "  +2.17%  "; ForceText = $false }
    @{ Cell = "D26"; Value = "2.419"; ForceText = $true }
    @{ Cell = "E26"; Value = "  +0.99%  "; ForceText = $false }
    @{ Cell = "E27"; Value = "  +0.65%  "; ForceText = $false }
    @{ Cell = "D28"; Value = "151.47"; ForceText = $true }
    @{ Cell = "E28"; Value = "  +1.06%  "; ForceText = $false }
    @{ Cell = "D29"; Value = "5.367"; ForceText = $true }
    @{ Cell = "E29"; Value = "  +2.73%  "; ForceText = $false }
    @{ Cell = "D30"; Value = "135.94"; ForceText = $true }
    @{ Cell = "E30"; Value = "  +1.54%  "; ForceText = $false }
    @{ Cell = "D31"; Value = "2.392"; ForceText = $true }
    @{ Cell = "E31"; Value = "  -1.25%  "; ForceText = $false }
    @{ Cell = "D32"; Value = "1.825.97"; ForceText = $true }
    @{ Cell = "E32"; Value = "  +2.51%  "; ForceText = $false }
    @{ Cell = "E33"; Value = "  +0.11%  "; ForceText = $false }
    @{ Cell = "D34"; Value = "0.9787"; ForceText = $true }
    @{ Cell = "E34"; Value = "  -0.73%  "; ForceText = $false }
    @{ Cell = "E36"; Value = "  +0.44%  "; ForceText = $false }
    @{ Cell = "D37"; Value = "0.07383"; ForceText = $true }
    @{ Cell = "E37"; Value = "  -2.73%  "; ForceText = $false }
    @{ Cell = "D38"; Value = "0.2557"; ForceText = $true }
    @{ Cell = "E38"; Value = "  +1.96%  "; ForceText = $false }
    @{ Cell = "D39"; Value = "6.213"; ForceText = $true }
    @{ Cell = "E39"; Value = "  +0.89%  "; ForceText = $false }
    @{ Cell = "D40"; Value = "0.08891"; ForceText = $true }
    @{ Cell = "E40"; Value = "  +1.11%  "; ForceText = $false }
    @{ Cell = "D41"; Value = "1.385"; ForceText = $true }
    @{ Cell = "E41"; Value = "  +1.73%  "; ForceText = $false }
    @{ Cell = "D42"; Value = "0.7140"; ForceText = $true }
    @{ Cell = "B43"; Value = "Aptos"; ForceText = $false }
    @{ Cell = "C43"; Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; ForceText = $false }
    @{ Cell = "D43"; Value = "12.57"; ForceText = $true }
    @{ Cell = "E43"; Value = "  +1.21%  "; ForceText = $false }
    @{ Cell = "B44"; Value = "EnergySwap"; ForceText = $false }
    @{ Cell = "C44"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; ForceText = $false }
    @{ Cell = "D44"; Value = "16.30"; ForceText = $true }
    @{ Cell = "E44"; Value = "  +5.57%  "; ForceText = $false }
    @{ Cell = "D45"; Value = "0.6573"; ForceText = $true }
    @{ Cell = "E45"; Value = "  -0.07%  "; ForceText = $false }
    @{ Cell = "E46"; Value = "  +1.59%  "; ForceText = $false }
    @{ Cell = "D47"; Value = "4.044"; ForceText = $true }
    @{ Cell = "E47"; Value = "  +0.68%  "; ForceText = $false }
    @{ Cell = "D48"; Value = "0.9992"; ForceText = $true }
    @{ Cell = "E48"; Value = "  +0.13%  "; ForceText = $false }
    @{ Cell = "D49"; Value = "0.07999"; ForceText = $true }
    @{ Cell = "E49"; Value = "  +0.51%  "; ForceText = $false }
    @{ Cell = "D50"; Value = "130.29"; ForceText = $true }
    @{ Cell = "E50"; Value = "  -1.54%  "; ForceText = $false }
    @{ Cell = "D51"; Value = "1.211"; ForceText = $true }
    @{ Cell = "E51"; Value = "  +0.13%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Apply a text number format so the numeric-looking string is not
        # reinterpreted as a number, then restore the original (default) style
        # so no stray formatting is left behind on the cell.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
